$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at row 4, shifting existing data rows (old 4..63) down to (5..64)
$ws.Rows("4:4").Insert()

# Populate the new row 4 with this week's data (same market/category metadata,
# new date + new prices for Jengibre at Vega Monumental Concepción)
$ws.Cells.Item(4, 1).Value() = 11
$ws.Cells.Item(4, 2).Value() = "Vega Monumental Concepción"
$ws.Cells.Item(4, 3).Value() = "Bíobío"
$ws.Cells.Item(4, 4).Value() = "8/22/2023"
$ws.Cells.Item(4, 5).Value() = 8
$ws.Cells.Item(4, 6).Value() = 100114007
$ws.Cells.Item(4, 7).Value() = "Jengibre"
$ws.Cells.Item(4, 8).Value() = "Sin especificar"
$ws.Cells.Item(4, 9).Value() = "Primera"
$ws.Cells.Item(4, 10).Value() = 40
$ws.Cells.Item(4, 11).Value() = 17500
$ws.Cells.Item(4, 12).Value() = 18000
$ws.Cells.Item(4, 13).Value() = 17750
$ws.Cells.Item(4, 14).Value() = "$/caja 13 kilos"
$ws.Cells.Item(4, 15).Value() = "Perú"
$ws.Cells.Item(4, 16).Value() = 1365
$ws.Cells.Item(4, 17).Value() = 13
$ws.Cells.Item(4, 18).Value() = "Hortaliza"
